$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "Bitcoin"
$ws.Cells.Item(2, 3).Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Cells.Item(2, 4).Value = "25.933.47"
$ws.Cells.Item(2, 5).Value = "  -0.29%  "

$ws.Cells.Item(3, 2).Value = "Ethereum"
$ws.Cells.Item(3, 3).Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Cells.Item(3, 4).Value = "1.637.61"
$ws.Cells.Item(3, 5).Value = "  +0.18%  "

$ws.Cells.Item(4, 2).Value = "TetherUSD"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.93%  "

$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "214.58"
$ws.Range("D5").Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.07%  "

$ws.Cells.Item(6, 2).Value = "XRP"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.506"
$ws.Range("D6").Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.50%  "

$ws.Cells.Item(7, 2).Value = "USDC"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.89%  "

$ws.Cells.Item(8, 2).Value = "Cardano"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.254"
$ws.Range("D8").Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.88%  "

$ws.Cells.Item(9, 2).Value = "Dogecoin"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.0636"
$ws.Range("D9").Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.68%  "

$ws.Cells.Item(10, 2).Value = "Solana"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.58"
$ws.Range("D10").Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.61%  "

$ws.Cells.Item(11, 2).Value = "TRON"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.91%  "

$ws.Cells.Item(12, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(12, 4).Value = "1.864.52"
$ws.Cells.Item(12, 5).Value = "  +0.17%  "

$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.24"
$ws.Range("D13").Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.07%  "

$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value = "1.638.88"
$ws.Cells.Item(14, 5).Value = "  +1.79%  "

$ws.Cells.Item(15, 2).Value = "Polygon"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.543"
$ws.Range("D15").Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -1.54%  "

$ws.Cells.Item(16, 2).Value = "ShibaInu"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(16, 4).Value = "0.0₃0757"
$ws.Cells.Item(16, 5).Value = "  -0.81%  "

$ws.Cells.Item(17, 2).Value = "Litecoin"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "62.50"
$ws.Range("D17").Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.03%  "

$ws.Cells.Item(18, 2).Value = "WrappedBTC"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(18, 4).Value = "25.956.47"
$ws.Cells.Item(18, 5).Value = "  -0.18%  "

$ws.Cells.Item(19, 2).Value = "Dai"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.00"
$ws.Range("D19").Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.91%  "

$ws.Cells.Item(20, 2).Value = "BitcoinCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "193.43"
$ws.Range("D20").Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.56%  "

$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "4.38"
$ws.Range("D21").Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.72%  "

$ws.Cells.Item(22, 2).Value = "Avalanche"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "9.91"
$ws.Range("D22").Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.76%  "

$ws.Cells.Item(23, 2).Value = "Chainlink"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.26"
$ws.Range("D23").Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -1.72%  "

$ws.Cells.Item(24, 2).Value = "Toncoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.80"
$ws.Range("D24").Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.21%  "

$ws.Cells.Item(25, 2).Value = "Monero"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "144.01"
$ws.Range("D25").Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.80%  "

$ws.Cells.Item(26, 2).Value = "BinanceUSD"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.88%  "

$ws.Cells.Item(27, 2).Value = "Stellar"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.126"
$ws.Range("D27").Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.63%  "

$ws.Cells.Item(28, 2).Value = "Cosmos"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.84"
$ws.Range("D28").Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.40%  "

$ws.Cells.Item(29, 2).Value = "EthereumClassic"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.45"
$ws.Range("D29").Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -1.09%  "

$ws.Cells.Item(30, 2).Value = "PancakeSwap"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.24"
$ws.Range("D30").Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.05%  "

$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.0499"
$ws.Range("D31").Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.26%  "

$ws.Cells.Item(32, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.29"
$ws.Range("D32").Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.30%  "

$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.21"
$ws.Range("D33").Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.76%  "

$ws.Cells.Item(34, 2).Value = "LidoDAOToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.55"
$ws.Range("D34").Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -2.99%  "

$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.18%  "

$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.902"
$ws.Range("D36").Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.49%  "

$ws.Cells.Item(37, 2).Value = "Maker"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(37, 4).Value = "1.135.96"
$ws.Cells.Item(37, 5).Value = "  -0.73%  "

$ws.Cells.Item(38, 2).Value = "ImmutableX"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.544"
$ws.Range("D38").Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.10%  "

$ws.Cells.Item(39, 2).Value = "MXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.46"
$ws.Range("D39").Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.05%  "

$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0157"
$ws.Range("D40").Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.40%  "

$ws.Cells.Item(41, 2).Value = "Quant"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "99.39"
$ws.Range("D41").Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.94%  "

$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.797"
$ws.Range("D42").Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.25%  "

$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.41"
$ws.Range("D43").Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -2.67%  "

$ws.Cells.Item(44, 2).Value = "RocketPoolETH"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(44, 4).Value = "1.773.84"
$ws.Cells.Item(44, 5).Value = "  +0.06%  "

$ws.Cells.Item(45, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(45, 4).Value = "0.0₆0115"
$ws.Cells.Item(45, 5).Value = "  +4.78%  "

$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "56.46"
$ws.Range("D46").Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.62%  "

$ws.Cells.Item(47, 2).Value = "Cronos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0529"
$ws.Range("D47").Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.32%  "

$ws.Cells.Item(48, 2).Value = "RenderToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.45"
$ws.Range("D48").Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.36%  "

$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.415"
$ws.Range("D49").Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.00%  "

$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.64"
$ws.Range("D50").Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.60%  "

$ws.Cells.Item(51, 2).Value = "Algorand"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0961"
$ws.Range("D51").Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.69%  "
